$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns I and J mirror the existing header style (bold/border/center),
# so copy H1's formatting onto I1:J1 before writing the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-5.
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 7
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 8
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 8
